$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.151.83"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "2.027.74"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.06"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.39"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  -4.50%  "
$ws.Range("D12").Value = "2.327.94"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.20"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -3.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.28"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.746"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.19"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D17").Value = "1.986.34"
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("D18").Value = "37.098.60"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.53"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +7.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.88"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.48"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.53"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.23"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -4.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.127"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.70"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -3.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.117"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.53"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.49"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.59"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +6.22%  "
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("D40").Value = "1.467.71"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.63"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.49"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0912"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "2.211.96"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.62"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -8.41%  "
